$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The round-trip loader maps an empty shared-string cell (F1) onto shared
# string index 0 ("number"); explicitly clear it so it stays blank as in
# the source file.
$ws.Cells.Item(1, 6).Value = ""

# Add "NA" values in column E (duplicate_image_filename) for data rows 2
# through 21, per the commit "add the NA's under duplicate_image_filename".
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
